# Mark the first 11 field records (rows 2-12) in the layout sheet as
# "Obrigatorio" (required) by changing column E from "N" to "S".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E12").Value = "S"
